# Apply updated crypto price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "41.165.24"; E = "  -3.53%  " },
    @{ Row = 3; D = "2.461.68"; E = "  -2.69%  " },
    @{ Row = 4; D = $null; E = "  +0.04%  " },
    @{ Row = 5; D = "311.71"; E = "  +0.71%  " },
    @{ Row = 6; D = "93.74"; E = "  -6.22%  " },
    @{ Row = 7; D = "0.551"; E = "  -2.90%  " },
    @{ Row = 8; D = $null; E = "  +0.02%  " },
    @{ Row = 9; D = $null; E = "  -4.78%  " },
    @{ Row = 10; D = $null; E = "  -6.99%  " },
    @{ Row = 11; D = $null; E = "  -3.27%  " },
    @{ Row = 12; D = "0.108"; E = "  -1.34%  " },
    @{ Row = 13; D = "6.95"; E = "  -5.08%  " },
    @{ Row = 14; D = "2.840.76"; E = "  -2.72%  " },
    @{ Row = 15; D = "2.459.62"; E = "  -4.28%  " },
    @{ Row = 16; D = "14.79"; E = "  -3.06%  " },
    @{ Row = 17; D = "0.782"; E = "  -3.46%  " },
    @{ Row = 18; D = "41.125.55"; E = "  -3.61%  " },
    @{ Row = 19; D = $null; E = "  -6.44%  " },
    @{ Row = 20; D = $null; E = "  -3.06%  " },
    @{ Row = 21; D = "11.13"; E = "  -8.84%  " },
    @{ Row = 22; D = "68.34"; E = "  -1.37%  " },
    @{ Row = 23; D = "235.11"; E = "  -3.23%  " },
    @{ Row = 24; D = $null; E = "  -4.26%  " },
    @{ Row = 25; D = $null; E = "  -0.11%  " },
    @{ Row = 26; D = $null; E = "  -6.02%  " },
    @{ Row = 27; D = $null; E = "  -5.87%  " },
    @{ Row = 28; D = $null; E = "  -6.39%  " },
    @{ Row = 29; D = $null; E = "  -5.68%  " },
    @{ Row = 30; D = "36.22"; E = "  -5.80%  " },
    @{ Row = 31; D = "152.43"; E = "  -4.17%  " },
    @{ Row = 32; D = "5.47"; E = "  -4.91%  " },
    @{ Row = 33; D = $null; E = "  -5.46%  " },
    @{ Row = 34; D = $null; E = "  -3.23%  " },
    @{ Row = 35; D = $null; E = "  -5.01%  " },
    @{ Row = 36; D = $null; E = "  -3.63%  " },
    @{ Row = 37; D = $null; E = "  -4.17%  " },
    @{ Row = 38; D = "16.93"; E = "  -7.74%  " },
    @{ Row = 39; D = $null; E = "  -2.97%  " },
    @{ Row = 40; D = $null; E = "  -8.23%  " },
    @{ Row = 41; D = "4.18"; E = "  -1.90%  " },
    @{ Row = 42; D = $null; E = "  +0.15%  " },
    @{ Row = 43; D = "20.03"; E = "  -10.50%  " },
    @{ Row = 44; D = "1.972.14"; E = "  -1.15%  " },
    @{ Row = 45; D = $null; E = "  -5.18%  " },
    @{ Row = 46; D = "3.02"; E = "  -7.87%  " },
    @{ Row = 47; D = $null; E = "  -2.23%  " },
    @{ Row = 48; D = "68.93"; E = "  -4.27%  " },
    @{ Row = 49; D = "96.69"; E = "  -3.98%  " },
    @{ Row = 50; D = $null; E = "  -6.51%  " },
    @{ Row = 51; D = "73.72"; E = "  -6.96%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel stores the price as text rather
        # than silently re-interpreting it as a number (these values use
        # "." as both decimal and thousands separators, e.g. "41.165.24").
        $cell = $ws.Cells.Item($row, 4)
        $cell.Value = "'" + $u.D
        $cell.Style = "Normal"
    }

    if ($null -ne $u.E) {
        $ws.Cells.Item($row, 5).Value = $u.E
    }
}
